$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above row 416, pushing the existing 416-503 block
# down to 419-506 (matches Excel's own row-shift semantics / picks up the
# D-column date style from the row above, same as the rest of the sheet).
$ws.Rows("416:418").Insert()

# --- New row 416 ---
$ws.Range("A416").Value = 10
$ws.Range("B416").Value = "Vega Modelo de Temuco"
$ws.Range("C416").Value = "La Araucanía"
$ws.Range("D416").Value = 44511
$ws.Range("E416").Value = 9
$ws.Range("F416").Value = 100114001
$ws.Range("G416").Value = "Papa"
$ws.Range("H416").Value = "Asterix"
$ws.Range("I416").Value = "1a (guarda)"
$ws.Range("J416").Value = 1000
$ws.Range("K416").Value = 7000
$ws.Range("L416").Value = 8000
$ws.Range("M416").Value = 7500
$ws.Range("N416").Value = "$/malla 25 kilos"
$ws.Range("O416").Value = "Provincia de Cautín"
$ws.Range("P416").Value = 300
$ws.Range("Q416").Value = 25
$ws.Range("R416").Value = "Hortaliza"

# --- New row 417 ---
$ws.Range("A417").Value = 10
$ws.Range("B417").Value = "Vega Modelo de Temuco"
$ws.Range("C417").Value = "La Araucanía"
$ws.Range("D417").Value = 44511
$ws.Range("E417").Value = 9
$ws.Range("F417").Value = 100114001
$ws.Range("G417").Value = "Papa"
$ws.Range("H417").Value = "Pehuenche"
$ws.Range("I417").Value = "1a nueva(o)"
$ws.Range("J417").Value = 600
$ws.Range("K417").Value = 13000
$ws.Range("L417").Value = 13000
$ws.Range("M417").Value = 13000
$ws.Range("N417").Value = "$/saco 25 kilos"
$ws.Range("O417").Value = "Provincia de Cautín"
$ws.Range("P417").Value = 520
$ws.Range("Q417").Value = 25
$ws.Range("R417").Value = "Hortaliza"

# --- New row 418 ---
$ws.Range("A418").Value = 10
$ws.Range("B418").Value = "Vega Modelo de Temuco"
$ws.Range("C418").Value = "La Araucanía"
$ws.Range("D418").Value = 44511
$ws.Range("E418").Value = 9
$ws.Range("F418").Value = 100114001
$ws.Range("G418").Value = "Papa"
$ws.Range("H418").Value = "Rosara"
$ws.Range("I418").Value = "1a (guarda)"
$ws.Range("J418").Value = 500
$ws.Range("K418").Value = 7000
$ws.Range("L418").Value = 8000
$ws.Range("M418").Value = 7600
$ws.Range("N418").Value = "$/malla 25 kilos"
$ws.Range("O418").Value = "Provincia de Cautín"
$ws.Range("P418").Value = 304
$ws.Range("Q418").Value = 25
$ws.Range("R418").Value = "Hortaliza"
